$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings are appended in this order: JavaFX, Add new Rating, Android

# Row 29: fill in the C/D/E/F/G columns (JavaFX Appl Verbessern block)
$ws.Range("C29").Value = "JavaFX Appl Verbessern"
$ws.Range("D29").Value2 = 42753
$ws.Range("E29").Value2 = 0.59027777777777779
$ws.Range("F29").Value2 = 0.65972222222222221
$ws.Range("G29").Value2 = 1

# Row 29: fill in the M/N/O/P/Q columns (Add new Rating einbauen block)
$ws.Range("M29").Value = "Add new Rating einbauen"
$ws.Range("N29").Value2 = 42753
$ws.Range("O29").Value2 = 0.59027777777777779
$ws.Range("P29").Value2 = 0.65972222222222221
$ws.Range("Q29").Value2 = 1

# Row 24: fill in the W/X/Y/Z/AA columns (Android Appl Verbesserung block)
$ws.Range("W24").Value = "Android Appl Verbesserung"
$ws.Range("X24").Value2 = 42753
$ws.Range("Y24").Value2 = 0.59027777777777779
$ws.Range("Z24").Value2 = 0.65972222222222221
$ws.Range("AA24").Value2 = 1

# Update sheet view: selection state changed in the saved workbook
$ws.Range("W25").Select()
